$wb = $excel.ActiveWorkbook

# --- test_case_1: swap the two course rows (row 2 <-> row 3) ---
$ws1 = $wb.Worksheets.Item("test_case_1")
$ws1.Range("A2").Value = "IBM Full Stack Software Developer"
$ws1.Range("B2").Value = 4.6
$ws1.Range("A3").Value = "Meta Front-End Developer"
$ws1.Range("B3").Value = 4.7

# --- test_case_2: refresh the language-learning pivot counts ---
$ws2 = $wb.Worksheets.Item("test_case_2")

$ws2.Range("C5").Value = "122"
$ws2.Range("C7").Value = "733"
$ws2.Range("C8").Value = "684"
$ws2.Range("C10").Value = "124"
$ws2.Range("C12").Value = "0"
$ws2.Range("C13").Value = "819"
$ws2.Range("C14").Value = "496"
$ws2.Range("C15").Value = "583"
$ws2.Range("A17").Value = "Chinese (Traditional)"
$ws2.Range("C17").Value = "5"
$ws2.Range("C18").Value = "4"
$ws2.Range("B19").Value = "Mixed"
$ws2.Range("C19").Value = "3"
$ws2.Range("B20").Value = "Mixed(120)"
$ws2.Range("C20").Value = "120"
$ws2.Range("A21").Value = "Dutch"
$ws2.Range("C21").Value = "907"
$ws2.Range("A22").Value = "Spanish (1,062)"
$ws2.Range("B22").Value = "Intermediate"
$ws2.Range("C22").Value = "754"
$ws2.Range("B23").Value = "Advanced"
$ws2.Range("C23").Value = "452"
$ws2.Range("C24").Value = "537"
$ws2.Range("B25").Value = "Mixed(131)"
$ws2.Range("C25").Value = "131"
